$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.019999999999999
$ws.Range("C2").Value = 1.045888919939244
$ws.Range("D2").Value = 1.051881940192188
$ws.Range("E2").Value = 1.053268210991738
$ws.Range("F2").Value = 1.063311325738323
$ws.Range("I2").Value = 1.044905204617404
$ws.Range("J2").Value = 1.050946246141702
$ws.Range("K2").Value = 1.054632212339489
$ws.Range("L2").Value = 1.056014648806958
$ws.Range("M2").Value = 1.066030308212695
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.046711308515427
$ws.Range("D3").Value = 1.052508943949809
$ws.Range("E3").Value = 1.05397996817924
$ws.Range("F3").Value = 1.06405630977827
$ws.Range("I3").Value = 1.045091769524372
$ws.Range("J3").Value = 1.051417107097217
$ws.Range("K3").Value = 1.055072428311693
$ws.Range("L3").Value = 1.056539676370492
$ws.Range("M3").Value = 1.066590454118424
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.047244254333979
$ws.Range("D4").Value = 1.052915360384767
$ws.Range("E4").Value = 1.054441585296461
$ws.Range("F4").Value = 1.06453943265089
$ws.Range("I4").Value = 1.045211735609936
$ws.Range("J4").Value = 1.051721900334248
$ws.Range("K4").Value = 1.055357287291029
$ws.Range("L4").Value = 1.056879790596832
$ws.Range("M4").Value = 1.066953303789797
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.047468495697134
$ws.Range("D5").Value = 1.053086384250075
$ws.Range("E5").Value = 1.054635901507858
$ws.Range("F5").Value = 1.064742791170016
$ws.Range("I5").Value = 1.045261988085564
$ws.Range("J5").Value = 1.051850061463764
$ws.Range("K5").Value = 1.05547704294549
$ws.Range("L5").Value = 1.057022865531828
$ws.Range("M5").Value = 1.067105939343362
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.047506157951058
$ws.Range("D6").Value = 1.053115109611215
$ws.Range("E6").Value = 1.054668542795457
$ws.Range("F6").Value = 1.064776950773773
$ws.Range("I6").Value = 1.045270415044979
$ws.Range("J6").Value = 1.051871581768095
$ws.Range("K6").Value = 1.055497150460956
$ws.Range("L6").Value = 1.057046893728762
$ws.Range("M6").Value = 1.067131572961319
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.047247249910131
$ws.Range("D7").Value = 1.052917644962622
$ws.Range("E7").Value = 1.054444180770066
$ws.Range("F7").Value = 1.064542148943306
$ws.Range("I7").Value = 1.045212407799208
$ws.Range("J7").Value = 1.051723612728933
$ws.Range("K7").Value = 1.055358887470022
$ws.Range("L7").Value = 1.056881702015849
$ws.Range("M7").Value = 1.066955342948058
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.046166682111975
$ws.Range("D8").Value = 1.052093692315917
$ws.Range("E8").Value = 1.053508531429312
$ws.Range("F8").Value = 1.063562874098921
$ws.Range("I8").Value = 1.04496841072533
$ws.Range("J8").Value = 1.05110535120631
$ws.Range("K8").Value = 1.054780982510441
$ws.Range("L8").Value = 1.056192003407663
$ws.Range("M8").Value = 1.066219528760589
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.044268835320233
$ws.Range("D9").Value = 1.050647252196431
$ws.Range("E9").Value = 1.051868025452026
$ws.Range("F9").Value = 1.061845542554791
$ws.Range("I9").Value = 1.044532717664342
$ws.Range("J9").Value = 1.050016836369524
$ws.Range("K9").Value = 1.053762777592207
$ws.Range("L9").Value = 1.05497969397541
$ws.Range("M9").Value = 1.064926051517513
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.043007919792841
$ws.Range("D10").Value = 1.049686751307041
$ws.Range("E10").Value = 1.050780007453413
$ws.Range("F10").Value = 1.060706345338235
$ws.Range("I10").Value = 1.04423844411677
$ws.Range("J10").Value = 1.049291876821326
$ws.Range("K10").Value = 1.053084151040742
$ws.Range("L10").Value = 1.054173616833449
$ws.Range("M10").Value = 1.064065932882404
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.042462976828691
$ws.Range("D11").Value = 1.049271767751481
$ws.Range("E11").Value = 1.050310250015834
$ws.Range("F11").Value = 1.060214436864409
$ws.Range("I11").Value = 1.044110126356056
$ws.Range("J11").Value = 1.048978149267948
$ws.Range("K11").Value = 1.052790358035929
$ws.Range("L11").Value = 1.053825101364987
$ws.Range("M11").Value = 1.063694035108966
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.042260719187369
$ws.Range("D12").Value = 1.049117764540866
$ws.Range("E12").Value = 1.050135967836318
$ws.Range("F12").Value = 1.060031928385335
$ws.Range("I12").Value = 1.044062329703447
$ws.Range("J12").Value = 1.048861646084572
$ws.Range("K12").Value = 1.05268124020926
$ws.Range("L12").Value = 1.053695727126297
$ws.Range("M12").Value = 1.063555978475601
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.042304096937394
$ws.Range("D13").Value = 1.049150792377916
$ws.Range("E13").Value = 1.050173342570563
$ws.Range("F13").Value = 1.060071067615682
$ws.Range("I13").Value = 1.044072588286116
$ws.Range("J13").Value = 1.048886635068348
$ws.Range("K13").Value = 1.0527046458788
$ws.Range("L13").Value = 1.053723474697077
$ws.Range("M13").Value = 1.0635855883192
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.04244625490572
$ws.Range("D14").Value = 1.049259034927354
$ws.Range("E14").Value = 1.05029583956884
$ws.Range("F14").Value = 1.060199346396121
$ws.Range("I14").Value = 1.044106178193134
$ws.Range("J14").Value = 1.048968518474454
$ws.Range("K14").Value = 1.052781338110452
$ws.Range("L14").Value = 1.053814405610319
$ws.Range("M14").Value = 1.063682621606287
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.04253386417412
$ws.Range("D15").Value = 1.049325745366404
$ws.Range("E15").Value = 1.050371341457921
$ws.Range("F15").Value = 1.060278410844398
$ws.Range("I15").Value = 1.04412685634773
$ws.Range("J15").Value = 1.049018973455843
$ws.Range("K15").Value = 1.052828592090794
$ws.Range("L15").Value = 1.053870441785946
$ws.Range("M15").Value = 1.063742418035252
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.043044107975872
$ws.Range("D16").Value = 1.049714311945134
$ws.Range("E16").Value = 1.050811212572182
$ws.Range("F16").Value = 1.060739020764969
$ws.Range("I16").Value = 1.044246941347436
$ws.Range("J16").Value = 1.049312701901087
$ws.Range("K16").Value = 1.053103650434529
$ws.Range("L16").Value = 1.054196757777077
$ws.Range("M16").Value = 1.06409062602311
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.043364450703231
$ws.Range("D17").Value = 1.049958297120639
$ws.Range("E17").Value = 1.051087498101285
$ws.Range("F17").Value = 1.061028317916761
$ws.Range("I17").Value = 1.044322028340298
$ws.Range("J17").Value = 1.049497000469781
$ws.Range("K17").Value = 1.053276203487494
$ws.Range("L17").Value = 1.054401587869435
$ws.Range("M17").Value = 1.064309193300607
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.043551401583834
$ws.Range("D18").Value = 1.050100698146964
$ws.Range("E18").Value = 1.051248782004223
$ws.Range("F18").Value = 1.061197192228252
$ws.Range("I18").Value = 1.044365738868854
$ws.Range("J18").Value = 1.049604516432074
$ws.Range("K18").Value = 1.053376856144381
$ws.Range("L18").Value = 1.054521111982223
$ws.Range("M18").Value = 1.064436731811018
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.043615163971864
$ws.Range("D19").Value = 1.050149268184733
$ws.Range("E19").Value = 1.051303797842029
$ws.Range("F19").Value = 1.061254796357369
$ws.Range("I19").Value = 1.044380628342895
$ws.Range("J19").Value = 1.049641179559234
$ws.Range("K19").Value = 1.053411176967863
$ws.Range("L19").Value = 1.05456187507917
$ws.Range("M19").Value = 1.064480227899995
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.043330070562862
$ws.Range("D20").Value = 1.049932110641561
$ws.Range("E20").Value = 1.051057841690111
$ws.Range("F20").Value = 1.06099726538276
$ws.Range("I20").Value = 1.04431398115215
$ws.Range("J20").Value = 1.049477225133702
$ws.Range("K20").Value = 1.053257689610974
$ws.Range("L20").Value = 1.05437960633649
$ws.Range("M20").Value = 1.064285737717849
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.042404388538084
$ws.Range("D21").Value = 1.049227156309381
$ws.Range("E21").Value = 1.050259761534689
$ws.Range("F21").Value = 1.060161565726817
$ws.Range("I21").Value = 1.044096290487725
$ws.Range("J21").Value = 1.048944405034598
$ws.Range("K21").Value = 1.052758753872714
$ws.Range("L21").Value = 1.053787626503534
$ws.Range("M21").Value = 1.063654045423233
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.041823293143252
$ws.Range("D22").Value = 1.04878473544088
$ws.Range("E22").Value = 1.049759173505504
$ws.Range("F22").Value = 1.059637334021668
$ws.Range("I22").Value = 1.043958646442871
$ws.Range("J22").Value = 1.048609569541534
$ws.Range("K22").Value = 1.052445111608126
$ws.Range("L22").Value = 1.05341588849915
$ws.Range("M22").Value = 1.06325735508056
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.042131255182538
$ws.Range("D23").Value = 1.049019193491943
$ws.Range("E23").Value = 1.050024430565224
$ws.Range("F23").Value = 1.059915124138207
$ws.Range("I23").Value = 1.04403168718726
$ws.Range("J23").Value = 1.048787055607521
$ws.Range("K23").Value = 1.052611373284598
$ws.Range("L23").Value = 1.053612909403806
$ws.Range("M23").Value = 1.06346760208885
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.043345605158488
$ws.Range("D24").Value = 1.049943942910779
$ws.Range("E24").Value = 1.051071241740586
$ws.Range("F24").Value = 1.061011296278125
$ws.Range("I24").Value = 1.04431761759729
$ws.Range("J24").Value = 1.04948616070228
$ws.Range("K24").Value = 1.053266055218421
$ws.Range("L24").Value = 1.054389538690195
$ws.Range("M24").Value = 1.064296336125887
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.044758721779053
$ws.Range("D25").Value = 1.051020531590414
$ws.Range("E25").Value = 1.052291148088256
$ws.Range("F25").Value = 1.062288520174853
$ws.Range("I25").Value = 1.044646030023766
$ws.Range("J25").Value = 1.050298123345534
$ws.Range("K25").Value = 1.054025982992345
$ws.Range("L25").Value = 1.055292736591284
$ws.Range("M25").Value = 1.065260066091853
